# CheckProblemDevice.xlsx - "problem management" report template rework
#
# The report's header row (row 7 on the "DeviceTermProb" sheet) is being
# repurposed from a "problem/branch" oriented layout to a device
# "event/terminal" oriented layout:
#   B: Branch Name   -> Date Time
#   C: Terminal ID   -> Serial No
#   D: Location      -> Terminal ID
#   E: Problem Name  -> Terminal Name
#   F: Remark        -> Location
#   G: Date Time     -> Event Name
#
# Column widths are widened to comfortably fit the new (generally longer)
# header captions and their data, and the active selection is left on the
# last header cell (G7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the report header columns (row 7) ---
$ws.Range("B7").Value = "Date Time"
$ws.Range("C7").Value = "Serial No"
$ws.Range("D7").Value = "Terminal ID"
$ws.Range("E7").Value = "Terminal Name"
$ws.Range("F7").Value = "Location"
$ws.Range("G7").Value = "Event Name"

# --- Widen columns B, E, F, G for the new headers/content ---
$ws.Columns.Item(2).ColumnWidth = 27.42857142857143   # B -> ~28.125
$ws.Columns.Item(5).ColumnWidth = 38.57142857142857   # E -> ~39.25
$ws.Columns.Item(6).ColumnWidth = 64.28571428571428   # F -> 65
$ws.Columns.Item(7).ColumnWidth = 42.85714285714286   # G -> ~43.625

# --- Leave the active cell/selection on the last header cell ---
$null = $ws.Range("G7").Select()
